$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row before current row 4 so the old row4 (251231) shifts to row5
$ws.Rows.Item(4).Insert()

# Swap rows 2 and 3: row2 (252277/45847...) and row3 (251849/45846...)
$ws.Range("A2").Value = 251849
$ws.Range("B2").Value = 45846.58333333334
$ws.Range("A3").Value = 252277
$ws.Range("B3").Value = 45847.58333333334

# Fill in the newly inserted row 4 with the new data
$ws.Range("A4").Value = 252284
$ws.Range("B4").Value = 45846.58333333334
$ws.Range("B4").NumberFormat = $ws.Range("B3").NumberFormat
